$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Mark Create Test Passed (B) and Read Test Passed (C) as TRUE for all
# Zone entries (rows 2-10, Z01-Z09) now that validation/checks were added.
$ws.Range("B2:C10").Value = $true

# Reflect the reviewer's scroll position / active cell at save time.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
[void]$ws.Range("C16").Select()
